$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed crypto price / 1h-volume data (+ a 3-row and a 3-row resort of the
# coin list) as pulled by the scheduled GitHub Actions scraper.
#
# Column D ("Price") stores numeric-looking values (e.g. "35.534.79",
# "0.690") as plain TEXT in the source data. Writing such a string straight
# to .Value lets Excel auto-coerce it to a Number (dropping things like
# trailing zeros), so each Price write: snapshots the cells current Style,
# forces a text NumberFormat, assigns the value, then restores the original
# Style -- the cell ends up as plain text with no formatting left behind.

function Set-TextValue($range, $value) {
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $savedStyle
}

Set-TextValue $ws.Range("D2") "35.534.79"
$ws.Range("E2").Value = "  +1.37%  "
Set-TextValue $ws.Range("D3") "1.902.35"
$ws.Range("E3").Value = "  +3.02%  "
$ws.Range("E4").Value = "  +0.73%  "
Set-TextValue $ws.Range("D5") "245.61"
$ws.Range("E5").Value = "  +5.15%  "
$ws.Range("E6").Value = "  +1.86%  "
$ws.Range("E7").Value = "  +0.61%  "
Set-TextValue $ws.Range("D8") "42.19"
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("E9").Value = "  +2.86%  "
$ws.Range("E10").Value = "  +1.34%  "
Set-TextValue $ws.Range("D11") "0.0998"
$ws.Range("E11").Value = "  +1.59%  "
Set-TextValue $ws.Range("D12") "2.180.64"
$ws.Range("E12").Value = "  +3.12%  "
Set-TextValue $ws.Range("D13") "12.40"
$ws.Range("E13").Value = "  +7.40%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D14") "1.911.54"
$ws.Range("E14").Value = "  +3.37%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D15") "0.690"
$ws.Range("E15").Value = "  +2.33%  "
Set-TextValue $ws.Range("D16") "4.82"
$ws.Range("E16").Value = "  +2.75%  "
Set-TextValue $ws.Range("D17") "35.523.11"
$ws.Range("E17").Value = "  +1.21%  "
Set-TextValue $ws.Range("D18") "71.84"
$ws.Range("E18").Value = "  +2.45%  "
$ws.Range("E19").Value = "  +2.16%  "
Set-TextValue $ws.Range("D20") "243.34"
$ws.Range("E20").Value = "  +1.12%  "
Set-TextValue $ws.Range("D21") "12.42"
$ws.Range("E21").Value = "  +1.82%  "
Set-TextValue $ws.Range("D22") "4.88"
$ws.Range("E22").Value = "  +2.19%  "
$ws.Range("E23").Value = "  +0.60%  "
Set-TextValue $ws.Range("D24") "2.29"
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("E25").Value = "  +28.56%  "
Set-TextValue $ws.Range("D26") "171.73"
$ws.Range("E26").Value = "  -0.04%  "
Set-TextValue $ws.Range("D27") "8.57"
$ws.Range("E27").Value = "  +8.50%  "
Set-TextValue $ws.Range("D28") "17.92"
$ws.Range("E28").Value = "  +2.34%  "
$ws.Range("E29").Value = "  +0.61%  "
Set-TextValue $ws.Range("D30") "0.978"
$ws.Range("E30").Value = "  +28.22%  "
Set-TextValue $ws.Range("D31") "4.08"
$ws.Range("E31").Value = "  +2.81%  "
Set-TextValue $ws.Range("D32") "0.0564"
$ws.Range("E32").Value = "  +1.68%  "
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("E34").Value = "  +4.63%  "
$ws.Range("E35").Value = "  +5.71%  "
$ws.Range("E36").Value = "  +2.53%  "
Set-TextValue $ws.Range("D37") "1.30"
$ws.Range("E37").Value = "  +4.30%  "
$ws.Range("E38").Value = "  +3.71%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D39") "0.0205"
$ws.Range("E39").Value = "  +2.25%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D40") "91.10"
$ws.Range("E40").Value = "  +1.14%  "
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue $ws.Range("D41") "51.62"
$ws.Range("E41").Value = "  +50.60%  "
Set-TextValue $ws.Range("D42") "1.352.08"
$ws.Range("E42").Value = "  +0.42%  "
Set-TextValue $ws.Range("D43") "15.46"
$ws.Range("E43").Value = "  +6.27%  "
$ws.Range("E44").Value = "  +11.51%  "
$ws.Range("E45").Value = "  +2.45%  "
Set-TextValue $ws.Range("D46") "12.72"
$ws.Range("E46").Value = "  +9.24%  "
$ws.Range("E47").Value = "  +1.87%  "
Set-TextValue $ws.Range("D48") "2.75"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("E49").Value = "  +5.04%  "
Set-TextValue $ws.Range("D50") "2.086.99"
$ws.Range("E50").Value = "  +2.80%  "
$ws.Range("E51").Value = "  +2.32%  "
